# Add "Wins" / "Losses" / "Ties" team-record columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the last existing header cell (bold, bordered,
# centered) onto the three new header cells so they match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Every player row gets the same team record: 79 wins, 83 losses, 0 ties.
$lastRow = 56
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 79   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 83   # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
